$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InsurancePremium")

# Header cells G1, M1, T1 previously had a dedicated "Text" number format (plus the
# yellow fill / border already shared with the rest of row 1). That extra format is
# no longer needed, so re-apply the plain fill+border format already used by the
# other header cells (e.g. A1) without introducing a new style entry.
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("T1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Product_startdate (column T) text values change from DD/MM/YYYY to MM/DD/YYYY
$ws.Range("T2").Value = "06/15/2023"
$ws.Range("T3").Value = "06/15/2023"
$ws.Range("T4").Value = "06/15/2023"

# Product_insurancesum (column U) changes from a number to the text "3000000"
$ws.Range("U2").NumberFormat = "@"
$ws.Range("U2").Value = "3000000"
$ws.Range("U3").NumberFormat = "@"
$ws.Range("U3").Value = "3000000"
$ws.Range("U4").NumberFormat = "@"
$ws.Range("U4").Value = "3000000"

# Remove the leftover style-only empty rows (7-11) below the data
$ws.Range("A7:AB11").Delete()

# Scroll the sheet view over to column Q and select T7 (matches saved view state)
$ws.Range("T7").Select()
$excel.ActiveWindow.ScrollColumn = 17
